# Add an "Avarage" column (H) with AVERAGE(D:F) formulas to the Class 1A sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("H1").Value = "Avarage"

# Fill H2:H23 with an AVERAGE formula of the Eng/Maths/Sci columns (D:F) for each row
$ws.Range("H2").Formula = "=AVERAGE(D2:F2)"
$ws.Range("H3:H23").Formula = "=AVERAGE(D3:F3)"

# Update the selection to the new column's first data cell
$ws.Range("H3").Select()
